# 4.2.2.1a workbook update:
#  1. Refresh the Kyrgyz title (A1) to the new wording.
#  2. Rename "urban settlements" / "countryside" (column C, every oblast
#     block) to "urban" / "rural".
#  3. Add the 2023 column (N) with header + data, copying number format /
#     font styling from the 2022 column (M) of the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Kyrgyz title in A1 -------------------------------------------------
$ws.Cells.Item(1, 1).Value2 = "4.2.2.1а Балдарды мектепке чейин билим берүү менен камтуу"

# --- 2. "urban settlements" -> "urban", "countryside" -> "rural" ----------
foreach ($r in @(6, 9, 12, 15, 18, 21, 24, 27)) {
    $ws.Cells.Item($r, 3).Value2 = "urban"
}
foreach ($r in @(7, 10, 13, 16, 19, 22, 25, 28)) {
    $ws.Cells.Item($r, 3).Value2 = "rural"
}

# --- 3. New column N: year 2023 -------------------------------------------
# Header cell (row 4) - same style as M4.
$ws.Cells.Item(4, 14).Value2 = 2023
$ws.Cells.Item(4, 13).Copy() | Out-Null
$ws.Cells.Item(4, 14).PasteSpecial(-4122) | Out-Null

# Data rows 5-29 - copy font/formatting from column M, then apply the 0.0
# number format used throughout the sheet (reuses/creates the matching
# cellXf automatically).
foreach ($pair in @(
    @(5, 28.34784779265912),
    @(6, 39.999446500300472),
    @(7, 23.198557483143556),
    @(8, 27.597876990321573),
    @(9, 47.175678010018999),
    @(10, 22.17579894112394),
    @(11, 24.100104034215697),
    @(12, 38.296287676015361),
    @(13, 19.410249509822766),
    @(14, 30.400174646089773),
    @(15, 44.562134629854725),
    @(16, 24.612036336109007),
    @(17, 39.266683582846994),
    @(18, 54.818496110630946),
    @(19, 36.591078066914498),
    @(20, 23.890520476423561),
    @(21, 16.93085228577992),
    @(22, 24.386979772654026),
    @(23, 28.919699950811605),
    @(24, 37.932834522359492),
    @(25, 26.985549456704376),
    @(26, 27.190143693828379),
    @(27, 54.006768771869439),
    @(28, 22.334624692306893),
    @(29, 36.01461582008131)
)) {
    $r = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($r, 14).Value2 = $val
    $ws.Cells.Item($r, 13).Copy() | Out-Null
    $ws.Cells.Item($r, 14).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 14).NumberFormat = "0.0"
}

# Last row (30) - already thick-bordered + numFmt164 on M30, copy as-is.
$ws.Cells.Item(30, 14).Value2 = 42.081208505725009
$ws.Cells.Item(30, 13).Copy() | Out-Null
$ws.Cells.Item(30, 14).PasteSpecial(-4122) | Out-Null
